$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: L4 reuses the exact same style as K4 (year header, plain integer) ---
$ws.Range("L4").Value = 2021
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- xf 14: rows [5] (font=font7) ---
$rowList = @(5)
$valMap = @{
  5 = 23.9
}
foreach ($r in $rowList) {
  $addr = "L" + $r
  $c = $ws.Range($addr)
  $c.Value = $valMap[$r]
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 9
  $c.Font.Bold = $true
  $c.Font.ColorIndex = 1
  $c.NumberFormat = "0.0"
  $c.HorizontalAlignment = -4152
  $c.WrapText = $false
}

# --- xf 15: rows [6] (font=font8) ---
$rowList = @(6)
$valMap = @{
  6 = 28.5
}
foreach ($r in $rowList) {
  $addr = "L" + $r
  $c = $ws.Range($addr)
  $c.Value = $valMap[$r]
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 9
  $c.Font.Bold = $false
  $c.NumberFormat = "0.0"
  $c.WrapText = $true
}

# --- xf 16: rows [7, 10, 12, 13, 15, 16, 18, 21, 22, 24, 25, 28] (font=font8) ---
$rowList = @(7,10,12,13,15,16,18,21,22,24,25,28)
$valMap = @{
  7 = 21.3
  10 = 23.3
  12 = 34.9
  13 = 19.2
  15 = 34.4
  16 = 26.7
  18 = 51.1
  21 = 10.8
  22 = 20.2
  24 = 56.8
  25 = 24.4
  28 = 18.6
}
foreach ($r in $rowList) {
  $addr = "L" + $r
  $c = $ws.Range($addr)
  $c.Value = $valMap[$r]
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 9
  $c.Font.Bold = $false
  $c.HorizontalAlignment = -4152
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
}

# --- xf 17: rows [8, 11, 14, 20, 23, 26] (font=font9) ---
$rowList = @(8,11,14,20,23,26)
$valMap = @{
  8 = 28.1
  11 = 23.2
  14 = 29.6
  20 = 19.5
  23 = 28.5
  26 = 21.7
}
foreach ($r in $rowList) {
  $addr = "L" + $r
  $c = $ws.Range($addr)
  $c.Value = $valMap[$r]
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 9
  $c.Font.Bold = $true
  $c.HorizontalAlignment = -4152
  $c.VerticalAlignment = -4108
  $c.WrapText = $true
}

# --- xf 18: rows [9, 19] (font=font8) ---
$rowList = @(9,19)
$valMap = @{
  9 = 44.5
  19 = 34.7
}
foreach ($r in $rowList) {
  $addr = "L" + $r
  $c = $ws.Range($addr)
  $c.Value = $valMap[$r]
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 9
  $c.Font.Bold = $false
  $c.NumberFormat = "0.0"
  $c.HorizontalAlignment = -4152
  $c.WrapText = $false
}

# --- xf 19: rows [17, 27] (font=font9) ---
$rowList = @(17,27)
$valMap = @{
  17 = 37.1
  27 = 35.1
}
foreach ($r in $rowList) {
  $addr = "L" + $r
  $c = $ws.Range($addr)
  $c.Value = $valMap[$r]
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 9
  $c.Font.Bold = $true
  $c.HorizontalAlignment = -4152
  $c.VerticalAlignment = -4108
  $c.WrapText = $false
}

# --- xf 20: rows [29] (font=font9) ---
$rowList = @(29)
$valMap = @{
  29 = 22.2
}
foreach ($r in $rowList) {
  $addr = "L" + $r
  $c = $ws.Range($addr)
  $c.Value = $valMap[$r]
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 9
  $c.Font.Bold = $true
  $c.NumberFormat = "0.0"
  $c.HorizontalAlignment = -4152
  $c.WrapText = $false
}

# --- xf 21: rows [30] (font=font9) ---
$rowList = @(30)
$valMap = @{
  30 = 29
}
foreach ($r in $rowList) {
  $addr = "L" + $r
  $c = $ws.Range($addr)
  $c.Value = $valMap[$r]
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 9
  $c.Font.Bold = $true
  $c.NumberFormat = "0.0"
  $c.HorizontalAlignment = -4152
  $c.WrapText = $true
  $c.Borders(-4107).LineStyle = 1
  $c.Borders(-4107).Weight = -4138
}

# --- Selection & view ---
$ws.Range("P19").Select()
